# Update plots for each sample
#
# The author lowered the minimum peak-height thresholds for three markers
# on the peak_table sheet (w_height / m_height, columns N/O). That, in
# turn, changes the downstream peak-detection results already baked into
# allele_table (per-allele peak stats), marker_table (per-marker
# genotype/phenotype calls) and genotype_result (the sample's overall
# genotype) -- this workbook stores pre-computed values rather than
# formulas, so each dependent sheet is updated explicitly here.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) peak_table - lower the detection thresholds (w_height / m_height)
# ---------------------------------------------------------------------
$peak = $wb.Worksheets.Item("peak_table")

$peak.Range("N2").Value = 900    # CYP2D6_14  w_height (wildtype G)
$peak.Range("O2").Value = 600    # CYP2D6_14  m_height (mutant A)

$peak.Range("N3").Value = 800    # CYP2D6_10B w_height (wildtype C)
# O3 (m_height) unchanged - stays 1000

$peak.Range("N4").Value = 600    # CYP2D6_49  w_height (wildtype T)
# O4 (m_height) unchanged - stays 1000

# ---------------------------------------------------------------------
# 2) allele_table - peaks that now clear the lowered threshold get
#    detected: min_height drops, is_detected flips to TRUE, and the
#    peak/size/height/status columns get filled in (message cleared).
# ---------------------------------------------------------------------
$allele = $wb.Worksheets.Item("allele_table")

# Row 2: CYP2D6_001 / CYP2D6_14, base G, wildtype
$allele.Range("K2").Value = 900
$allele.Range("M2").Value = $true
$allele.Range("N2").Value = 36
$allele.Range("O2").Value = 29.14
$allele.Range("P2").Value = 994
$allele.Range("Q2").Value = "ok"
$allele.Range("R2").Value = ""

# Row 3: CYP2D6_001 / CYP2D6_14, base A, mutant
$allele.Range("K3").Value = 600
$allele.Range("M3").Value = $true
$allele.Range("N3").Value = 35
$allele.Range("O3").Value = 30.91
$allele.Range("P3").Value = 705
$allele.Range("Q3").Value = "ok"
$allele.Range("R3").Value = ""

# Row 4: CYP2D6_002 / CYP2D6_10B, base C, wildtype
$allele.Range("K4").Value = 800
$allele.Range("M4").Value = $true
$allele.Range("N4").Value = 42
$allele.Range("O4").Value = 32.67
$allele.Range("P4").Value = 896
$allele.Range("Q4").Value = "ok"
$allele.Range("R4").Value = ""

# Row 6: CYP2D6_003 / CYP2D6_49, base T, wildtype
$allele.Range("K6").Value = 600
$allele.Range("M6").Value = $true
$allele.Range("N6").Value = 19
$allele.Range("O6").Value = 39.03
$allele.Range("P6").Value = 766
$allele.Range("Q6").Value = "ok"
$allele.Range("R6").Value = ""

# ---------------------------------------------------------------------
# 3) marker_table - genotype / phenotype calls for the three markers
#    that now have both alleles detected.
# ---------------------------------------------------------------------
$marker = $wb.Worksheets.Item("marker_table")

$marker.Range("G2").Value = "GA"
$marker.Range("H2").Value = "heterozygous"

$marker.Range("G3").Value = "CT"
$marker.Range("H3").Value = "heterozygous"

$marker.Range("G4").Value = "TT"
$marker.Range("H4").Value = "wildtype"

# ---------------------------------------------------------------------
# 4) genotype_result - overall sample genotype call
# ---------------------------------------------------------------------
$result = $wb.Worksheets.Item("genotype_result")
$result.Range("B2").Value = "*10B/*14B"
